$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell J1 (reuse formatting from I1), then set its value/text
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Q8"

# Row 2
$ws.Range("B2").Value = -0.9330598254158777
$ws.Range("C2").Value = 0.1669401745841223
$ws.Range("D2").Value = 0.3669401745841223
$ws.Range("E2").Value = -0.03305982541587771
$ws.Range("F2").Value = 0.2669401745841223
$ws.Range("G2").Value = -0.1330598254158777
$ws.Range("H2").Value = 0.6669401745841224

# Row 3
$ws.Range("B3").Value = -0.06208512593411045
$ws.Range("C3").Value = 0.1379148740658896
$ws.Range("D3").Value = -0.2620851259341104
$ws.Range("E3").Value = 0.03791487406588956
$ws.Range("F3").Value = -0.3620851259341105
$ws.Range("G3").Value = 0.4379148740658896

# Row 4
$ws.Range("B4").Value = 0.0543279172719193
$ws.Range("C4").Value = -0.3456720827280807
$ws.Range("D4").Value = -0.04567208272808071
$ws.Range("E4").Value = -0.4456720827280807
$ws.Range("F4").Value = 0.3543279172719193
$ws.Range("G4").Value = 0.0543279172719193
$ws.Range("H4").Value = -0.04567208272808071
$ws.Range("I4").Value = -0.4456720827280807
$ws.Range("J4").Value = -0.8456720827280808

# Row 5
$ws.Range("B5").Value = -0.8154625125417774
$ws.Range("C5").Value = -0.5154625125417773
$ws.Range("D5").Value = -0.9154625125417774
$ws.Range("E5").Value = -0.1154625125417773
$ws.Range("F5").Value = -0.4154625125417774
$ws.Range("G5").Value = -0.5154625125417773
$ws.Range("H5").Value = -0.9154625125417774
$ws.Range("I5").Value = -1.315462512541777

# Row 6
$ws.Range("B6").Value = -0.1813602613933202
$ws.Range("C6").Value = -0.5813602613933202
$ws.Range("D6").Value = 0.2186397386066798
$ws.Range("E6").Value = -0.0813602613933202
$ws.Range("F6").Value = -0.1813602613933202
$ws.Range("G6").Value = -0.5813602613933202
$ws.Range("H6").Value = -0.9813602613933202

# Row 7
$ws.Range("B7").Value = -0.4148081973238454
$ws.Range("C7").Value = 0.3851918026761547
$ws.Range("D7").Value = 0.08519180267615464
$ws.Range("E7").Value = -0.01480819732384536
$ws.Range("F7").Value = -0.4148081973238454
$ws.Range("G7").Value = -0.8148081973238454

# Row 8
$ws.Range("B8").Value = 0.4291840095081929
$ws.Range("C8").Value = 0.1291840095081928
$ws.Range("D8").Value = 0.02918400950819283
$ws.Range("E8").Value = -0.3708159904918072
$ws.Range("F8").Value = -0.7708159904918073
$ws.Range("G8").Value = -0.7708159904918073
$ws.Range("H8").Value = 0.9291840095081928
$ws.Range("I8").Value = -0.4708159904918072

# Row 9
$ws.Range("B9").Value = 0.06678455670716371
$ws.Range("C9").Value = -0.03321544329283629
$ws.Range("D9").Value = -0.4332154432928363
$ws.Range("E9").Value = -0.8332154432928364
$ws.Range("F9").Value = -0.8332154432928364
$ws.Range("G9").Value = 0.8667845567071637
$ws.Range("H9").Value = -0.5332154432928363

# Row 10
$ws.Range("B10").Value = 0.00001303303454188581
$ws.Range("C10").Value = -0.3999869669654582
$ws.Range("D10").Value = -0.7999869669654581
$ws.Range("E10").Value = -0.7999869669654581
$ws.Range("F10").Value = 0.9000130330345419
$ws.Range("G10").Value = -0.4999869669654581

# Row 11
$ws.Range("B11").Value = -0.406125572440377
$ws.Range("C11").Value = -0.8061255724403771
$ws.Range("D11").Value = -0.8061255724403771
$ws.Range("E11").Value = 0.8938744275596231
$ws.Range("F11").Value = -0.506125572440377

# Row 12
$ws.Range("B12").Value = -0.7512006267496926
$ws.Range("C12").Value = -0.7512006267496926
$ws.Range("D12").Value = 0.9487993732503075
$ws.Range("E12").Value = -0.4512006267496925

# Row 13
$ws.Range("B13").Value = -0.7522304086392605
$ws.Range("C13").Value = 0.9477695913607396
$ws.Range("D13").Value = -0.4522304086392604

# Row 14
$ws.Range("B14").Value = 1.371451352842971
$ws.Range("C14").Value = -0.02854864715702948

# Row 15
$ws.Range("B15").Value = -0.526054543893956

# Row 16 - no value changes (only the A16 label's underlying shared-string index shifts)
